$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the second row of the table (the "AC-AKB-00-000-000-P1B1-01" / Akbil / Yok / Yok record),
# shifting the rows below it upward and leaving the last row (13) blank.
$ws.Rows("2").Delete()

# Update the remembered selection to match the post-edit cursor position.
$ws.Range("I15").Select()
